$d = $word.ActiveDocument

$xmlFragment = @'
<w:p w14:paraId="05150551" w14:textId="473FDBAF" w:rsidR="00E75029" w:rsidRPr="00D23659" w:rsidRDefault="00EA637C" w:rsidP="00E75029"><w:pPr><w:jc w:val="center"/><w:rPr><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="091437D1" wp14:editId="1C927FD0"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>6350</wp:posOffset></wp:positionV><wp:extent cx="981075" cy="5457825"/><wp:effectExtent l="0" t="0" r="28575" b="28575"/><wp:wrapNone/><wp:docPr id="1152793138" name="Rectangle 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="981075" cy="5457825"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="15000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:fldSimple w:instr=" MERGEFIELD  image:mermaid()  \* MERGEFORMAT "><w:r><w:rPr><w:noProof/></w:rPr><w:t>«image:mermaid()»</w:t></w:r></w:fldSimple></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="091437D1" id="Rectangle 1" o:spid="_x0000_s1026" style="position:absolute;left:0;text-align:left;margin-left:0;margin-top:.5pt;width:77.25pt;height:429.75pt;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:center;mso-position-horizontal-relative:margin;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-height-percent:0;mso-height-relative:margin;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQCwgN7YYwIAAB4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVFFP2zAQfp+0/2D5fSSp6ICKFFUgpkkI&#10;0MrEs+vYJJLj885uk+7X7+ykKQK0h2kvztl39935y3e+vOpbw3YKfQO25MVJzpmyEqrGvpT859Pt&#10;l3POfBC2EgasKvleeX61/PzpsnMLNYMaTKWQEYj1i86VvA7BLbLMy1q1wp+AU5acGrAVgbb4klUo&#10;OkJvTTbL869ZB1g5BKm8p9ObwcmXCV9rJcOD1l4FZkpOvYW0Ylo3cc2Wl2LxgsLVjRzbEP/QRSsa&#10;S0UnqBsRBNti8w6qbSSCBx1OJLQZaN1Ile5AtynyN7dZ18KpdBcix7uJJv//YOX9bu0ekWjonF94&#10;MuMteo1t/FJ/rE9k7SeyVB+YpMOL8yI/m3MmyTU/nZ+dz+aRzeyY7dCHbwpaFo2SI/2MxJHY3fkw&#10;hB5CKO9YP1lhb1RswdgfSrOmooqzlJ2koa4Nsp2gnyqkVDYUg6sWlRqOi3mep79L/UwZqbsEGJF1&#10;Y8yEPQJE2b3HHnod42OqSsqakvO/NTYkTxmpMtgwJbeNBfwIwNCtxspD/IGkgZrIUug3PYVEcwPV&#10;/hEZwiBx7+RtQ7TfCR8eBZKmSf00p+GBFm2gKzmMFmc14O+PzmM8SY28nHU0IyX3v7YCFWfmuyUR&#10;XhSnp3Go0oYkMKMNvvZsXnvstr0G+mMFvQhOJjPGB3MwNUL7TOO8ilXJJayk2iWXAQ+b6zDMLj0I&#10;Uq1WKYwGyYlwZ9dORvBIcJTVU/8s0I3aC6TaezjMk1i8keAQGzMtrLYBdJP0eeR1pJ6GMGlofDDi&#10;lL/ep6jjs7b8AwAA//8DAFBLAwQUAAYACAAAACEAs2BKBtsAAAAGAQAADwAAAGRycy9kb3ducmV2&#10;LnhtbEyPT0/DMAzF70h8h8hI3FjKn06lNJ0mBIfd2Jh29hqvrUicqsm2wqfHO8HJen7Wez9Xi8k7&#10;daIx9oEN3M8yUMRNsD23Braf73cFqJiQLbrAZOCbIizq66sKSxvOvKbTJrVKQjiWaKBLaSi1jk1H&#10;HuMsDMTiHcLoMYkcW21HPEu4d/ohy+baY8/S0OFArx01X5ujN/CzOujsI74V2+XqOX/s1263Q2fM&#10;7c20fAGVaEp/x3DBF3SohWkfjmyjcgbkkSRbGRczf8pB7Q0U8ywHXVf6P379CwAA//8DAFBLAQIt&#10;ABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10u&#10;eG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5y&#10;ZWxzUEsBAi0AFAAGAAgAAAAhALCA3thjAgAAHgUAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9E&#10;b2MueG1sUEsBAi0AFAAGAAgAAAAhALNgSgbbAAAABgEAAA8AAAAAAAAAAAAAAAAAvQQAAGRycy9k&#10;b3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAADFBQAAAAA=&#10;" fillcolor="#5b9bd5 [3204]" strokecolor="#091723 [484]" strokeweight="1pt"><v:textbox><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:fldSimple w:instr=" MERGEFIELD  image:mermaid()  \* MERGEFORMAT "><w:r><w:rPr><w:noProof/></w:rPr><w:t>«image:mermaid()»</w:t></w:r></w:fldSimple></w:p></w:txbxContent></v:textbox><w10:wrap anchorx="margin"/></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="01Dautru"/><w:ind w:left="284"/></w:pPr></w:p>
'@

$marker = [char]0xAB + "image:mermaid()" + [char]0xBB

$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$marker*") {
        $p.Range.InsertXML($xmlFragment)
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not locate the image:mermaid() MERGEFIELD paragraph to replace."
}

Write-Output "Replaced image:mermaid() paragraph with floating textbox drawing: $found"
